# Update model and logic diagrams
#
# The ":AddressBookParser" lifeline label (two paragraphs: ":Address" and
# "BookParser") becomes a single paragraph made of two runs: ":" and
# "HealthBaseParser".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The label lives in shape #5 on the slide (a "Rectangle 62" shaped like a
# lifeline header) whose text currently reads ":AddressBookParser" across two
# paragraphs.
$shape = $s.Shapes.Item(5)
$tr = $shape.TextFrame.TextRange

# Collapse the two paragraphs into one, keeping the leading colon, then
# replace the class name portion with the new name. Characters() returns a
# sub-range so each edit becomes its own run instead of clobbering the whole
# text body.
$tr.Text = ":HealthBaseParser"
$colon = $tr.Characters(1, 1)
$colon.Text = ":"
$rest = $tr.Characters(2, 16)
$rest.Text = "HealthBaseParser"
